$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Table 1 (rows 6-8): "Win 11 - Ryzen 7950X ..." benchmark numbers
# ---------------------------------------------------------------------------
$ws.Cells.Item(6,4).Value = 856
$ws.Cells.Item(6,5).Value = 832
$ws.Cells.Item(6,6).Value = 1405
$ws.Cells.Item(6,7).Value = 1373

$ws.Cells.Item(7,5).Value = 830
$ws.Cells.Item(7,7).Value = 1400

# ---------------------------------------------------------------------------
# Table 2 (rows 14-16): "Mac M1 Pro ..." benchmark numbers
# ---------------------------------------------------------------------------
$ws.Cells.Item(14,4).Value = 2390
$ws.Cells.Item(14,5).Value = 2315
$ws.Cells.Item(14,6).Value = 3640
$ws.Cells.Item(14,7).Value = 3583

$ws.Cells.Item(15,4).Value = 2153
$ws.Cells.Item(15,5).Value = 2162
$ws.Cells.Item(15,6).Value = 3283
$ws.Cells.Item(15,7).Value = 3195

$ws.Cells.Item(16,4).Value = 2676
$ws.Cells.Item(16,5).Value = 2689
$ws.Cells.Item(16,6).Value = 4490
$ws.Cells.Item(16,7).Value = 4353

# ---------------------------------------------------------------------------
# Table 3 (rows 22-24): "Ubuntu 22.04 - Ryzen 7840HS ..." benchmark numbers
# ---------------------------------------------------------------------------
$ws.Cells.Item(22,4).Value = 1244
$ws.Cells.Item(22,4).Font.Bold = $true
$ws.Cells.Item(22,5).Value = 1154
$ws.Cells.Item(22,6).Value = 2851
$ws.Cells.Item(22,6).Font.Bold = $true
$ws.Cells.Item(22,7).Value = 2798

$ws.Cells.Item(23,4).Value = 1234
$ws.Cells.Item(23,5).Value = 1310
$ws.Cells.Item(23,5).Font.Bold = $false
$ws.Cells.Item(23,6).Value = 2832
$ws.Cells.Item(23,7).Value = 2792

$ws.Cells.Item(24,4).Value = 1419
$ws.Cells.Item(24,5).Value = 1444
$ws.Cells.Item(24,6).Value = 3755
$ws.Cells.Item(24,7).Value = 3804

# ---------------------------------------------------------------------------
# Table 4 (rows 30-32): "Win 11 - Ryzen 7840HS ..." benchmark numbers
# ---------------------------------------------------------------------------
$ws.Cells.Item(30,4).Value = 1206
$ws.Cells.Item(30,4).Font.Bold = $true
$ws.Cells.Item(30,5).Value = 1229
$ws.Cells.Item(30,6).Value = 2711
$ws.Cells.Item(30,7).Value = 2760

$ws.Cells.Item(31,4).Value = 1159
$ws.Cells.Item(31,5).Value = 1138
$ws.Cells.Item(31,6).Value = 2686
$ws.Cells.Item(31,7).Value = 2549

$ws.Cells.Item(32,4).Value = 2017
$ws.Cells.Item(32,5).Value = 2031
$ws.Cells.Item(32,6).Value = 4138
$ws.Cells.Item(32,7).Value = 4240

# ---------------------------------------------------------------------------
# Selection moved from K9 to G17 (as last recorded in the sheet view)
# ---------------------------------------------------------------------------
[void]$ws.Range("G17").Select()
